# Marksheet update: recompute student-answer columns, scoring row, and
# drop the now-unused third question-set / extra answer columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Scoring summary rows (9-12) -----------------------------------------
$ws.Range("A10").Style = "mtitleStyle"
$ws.Range("B10").Value = 20
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = 7
$ws.Range("E10").Value = 28

$ws.Range("A11").Style = "mtitleStyle"
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

$ws.Range("A12").Style = "mtitleStyle"
$ws.Range("B12").Value = 80
$ws.Range("C12").Value = -1
$ws.Range("E12").Value = "79/112"

# --- Drop the third Student/Correct Ans column pair (G:H), all rows ------
$ws.Range("G15:H40").Clear()

# --- Drop the second pair's rows below the first three questions ---------
$ws.Range("D19:E40").Clear()

# --- Fill in the first "Student Ans" column (A16:A40) ---------------------
# Mirrors the "Correct Ans" column (B) for every attempted question;
# unattempted rows are left blank (normalStyle), and the one missed
# question (row 35) is flagged with incorrectStyle.
function Set-StudentAns($row, $answer, $style) {
    $cell = $ws.Range("A$row")
    if ($answer) {
        $cell.Value = $answer
    }
    $cell.Style = $style
}

Set-StudentAns 16 "Option A" "correctStyle"
Set-StudentAns 17 $null      "normalStyle"
Set-StudentAns 18 "Option B" "correctStyle"
Set-StudentAns 19 "Option C" "correctStyle"
Set-StudentAns 20 $null      "normalStyle"
Set-StudentAns 21 "Option C" "correctStyle"
Set-StudentAns 22 "Option D" "correctStyle"
Set-StudentAns 23 "Option D" "correctStyle"
Set-StudentAns 24 $null      "normalStyle"
Set-StudentAns 25 "Option A" "correctStyle"
Set-StudentAns 26 $null      "normalStyle"
Set-StudentAns 27 "Option A" "correctStyle"
Set-StudentAns 28 "Option D" "correctStyle"
Set-StudentAns 29 $null      "normalStyle"
Set-StudentAns 30 "Option B" "correctStyle"
Set-StudentAns 31 $null      "normalStyle"
Set-StudentAns 32 "Option C" "correctStyle"
Set-StudentAns 33 "Option D" "correctStyle"
Set-StudentAns 34 "Option B" "correctStyle"
Set-StudentAns 35 "Option C" "incorrectStyle"
Set-StudentAns 36 "Option A" "correctStyle"
Set-StudentAns 37 $null      "normalStyle"
Set-StudentAns 38 "Option A" "correctStyle"
Set-StudentAns 39 "Option D" "correctStyle"
Set-StudentAns 40 "Option D" "correctStyle"

# --- Fill in the second "Student Ans" column for the 3 remaining rows ----
$ws.Range("D16").Value = "Option A"
$ws.Range("D16").Style = "correctStyle"
$ws.Range("D17").Value = "Option C"
$ws.Range("D17").Style = "correctStyle"
$ws.Range("D18").Value = "Option D"
$ws.Range("D18").Style = "correctStyle"
